$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "41.208.18"
$ws.Range("E2").Value = "  +2.90%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.248.25"
$ws.Range("E3").Value = "  +1.74%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'302.78"
$ws.Range("E5").Value = "  +3.15%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'91.19"
$ws.Range("E6").Value = "  +4.38%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.520"
$ws.Range("E7").Value = "  +2.18%  "

# Row 8 - USDC (price unchanged)
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 - Cardano (price unchanged)
$ws.Range("E9").Value = "  +2.31%  "

# Row 10 - OKB
$ws.Range("D10").Value = "'53.69"
$ws.Range("E10").Value = "  +8.24%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "'31.92"
$ws.Range("E11").Value = "  +6.75%  "

# Row 12 - Dogecoin
$ws.Range("D12").Value = "'0.0794"
$ws.Range("E12").Value = "  +2.11%  "

# Row 13 - TRON (price unchanged)
$ws.Range("E13").Value = "  +3.29%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'6.58"
$ws.Range("E14").Value = "  +1.91%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.598.18"
$ws.Range("E15").Value = "  +1.90%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "'14.10"
$ws.Range("E16").Value = "  +2.63%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.217.22"
$ws.Range("E17").Value = "  +0.04%  "

# Row 18 - Polygon (price unchanged)
$ws.Range("E18").Value = "  +3.23%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "41.141.50"
$ws.Range("E19").Value = "  +3.02%  "

# Row 20 - InternetComputer(DFINITY)
$ws.Range("D20").Value = "'12.06"
$ws.Range("E20").Value = "  +6.96%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0902"
$ws.Range("E21").Value = "  +2.03%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'5.86"
$ws.Range("E22").Value = "  +1.55%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'66.66"
$ws.Range("E23").Value = "  +2.16%  "

# Row 24 - BitcoinCash
$ws.Range("D24").Value = "'240.21"
$ws.Range("E24").Value = "  +0.94%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "'2.56"
$ws.Range("E25").Value = "  +4.78%  "

# Row 26 - Dai (price unchanged)
$ws.Range("E26").Value = "  -0.31%  "

# Row 27 - ImmutableX
$ws.Range("D27").Value = "'1.85"
$ws.Range("E27").Value = "  +2.75%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'23.73"
$ws.Range("E28").Value = "  +5.43%  "

# Row 29/30 - Toncoin and Cosmos swap places (Cosmos now row 29, Toncoin now row 30)
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'9.61"
$ws.Range("E29").Value = "  +4.74%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.10"
$ws.Range("E30").Value = "  -2.71%  "

# Row 31 - Monero
$ws.Range("D31").Value = "'158.41"
$ws.Range("E31").Value = "  +1.64%  "

# Row 32 - InjectiveProtocol
$ws.Range("D32").Value = "'33.32"
$ws.Range("E32").Value = "  +5.96%  "

# Row 33 - FirstDigitalUSD (price unchanged)
$ws.Range("E33").Value = "  -0.13%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "'5.16"
$ws.Range("E34").Value = "  +5.44%  "

# Row 35 - Hedera
$ws.Range("D35").Value = "'0.0733"
$ws.Range("E35").Value = "  +3.21%  "

# Row 36 - LidoDAOToken
$ws.Range("D36").Value = "'3.03"
$ws.Range("E36").Value = "  +7.15%  "

# Row 37 - WEMIXToken (price unchanged)
$ws.Range("E37").Value = "  +1.22%  "

# Row 38 - Celestia
$ws.Range("D38").Value = "'16.68"
$ws.Range("E38").Value = "  +8.74%  "

# Row 39 - Stellar (price unchanged)
$ws.Range("E39").Value = "  +2.66%  "

# Row 40 - Kaspa (price unchanged)
$ws.Range("E40").Value = "  +5.02%  "

# Row 41 - ARBITRUM
$ws.Range("D41").Value = "'1.78"
$ws.Range("E41").Value = "  +6.71%  "

# Row 42 - RenderToken (price unchanged)
$ws.Range("E42").Value = "  +4.58%  "

# Row 43/44 - Maker and EnergySwap swap places (EnergySwap now row 43, Maker now row 44)
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'20.32"
$ws.Range("E43").Value = "  +15.73%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.065.45"
$ws.Range("E44").Value = "  -2.97%  "

# Row 45 - VeChain (price unchanged)
$ws.Range("E45").Value = "  +3.46%  "

# Row 46 - FraxShare
$ws.Range("D46").Value = "'10.19"
$ws.Range("E46").Value = "  +5.80%  "

# Row 47 - NEARProtocol
$ws.Range("D47").Value = "'2.94"
$ws.Range("E47").Value = "  +11.26%  "

# Row 48 - ApeXProtocol (price unchanged)
$ws.Range("E48").Value = "  -4.05%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "2.467.88"
$ws.Range("E49").Value = "  +2.09%  "

# Row 50/51 - Stacks and TrustWalletToken swap places (TrustWalletToken now row 50, Stacks now row 51)
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "'1.14"
$ws.Range("E50").Value = "  +3.95%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.50"
$ws.Range("E51").Value = "  +1.30%  "
